# Auto-generated edit script
# Applies the 2022-05-16 daily crime data update across all affected sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("I2").Value = 2287
$ws.Range("H3").Value = 8349
$ws.Range("I3").Value = 2450
$ws.Range("E4").Value = 1963
$ws.Range("F4").Value = 1853
$ws.Range("G4").Value = 1428
$ws.Range("I4").Value = 594
$ws.Range("I5").Value = 217
$ws.Range("I6").Value = 2832
$ws.Range("E7").Value = 25967
$ws.Range("F7").Value = 24042
$ws.Range("G7").Value = 24651
$ws.Range("H7").Value = 25972
$ws.Range("I7").Value = 8380

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("I7").Value = 281
$ws.Range("I8").Value = 526
$ws.Range("I9").Value = 45
$ws.Range("I11").Value = 139
$ws.Range("I15").Value = 109
$ws.Range("I19").Value = 236
$ws.Range("G20").Value = 619
$ws.Range("I20").Value = 222
$ws.Range("I21").Value = 51
$ws.Range("I29").Value = 542
$ws.Range("I31").Value = 82
$ws.Range("I33").Value = 383
$ws.Range("I36").Value = 113
$ws.Range("H37").Value = 941
$ws.Range("I37").Value = 267
$ws.Range("I42").Value = 291
$ws.Range("I44").Value = 66
$ws.Range("I48").Value = 91
$ws.Range("I49").Value = 57
$ws.Range("I51").Value = 76
$ws.Range("I53").Value = 85
$ws.Range("I54").Value = 189
$ws.Range("E63").Value = 310
$ws.Range("F63").Value = 150
$ws.Range("I63").Value = 34
$ws.Range("I65").Value = 198
$ws.Range("I67").Value = 325
$ws.Range("I76").Value = 129
$ws.Range("I77").Value = 43
$ws.Range("I79").Value = 219
$ws.Range("I83").Value = 164
$ws.Range("I84").Value = 62
$ws.Range("I86").Value = 49
$ws.Range("I87").Value = 15
$ws.Range("I90").Value = 94
$ws.Range("I91").Value = 97
$ws.Range("I94").Value = 73
$ws.Range("I95").Value = 140
$ws.Range("I100").Value = 13
$ws.Range("E101").Value = 25967
$ws.Range("F101").Value = 24042
$ws.Range("G101").Value = 24651
$ws.Range("H101").Value = 25972
$ws.Range("I101").Value = 8380

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("I3").Value = 30
$ws.Range("I7").Value = 139

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("I2").Value = 165
$ws.Range("I3").Value = 141
$ws.Range("I4").Value = 33
$ws.Range("I7").Value = 526

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("I3").Value = 26
$ws.Range("I7").Value = 85

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("I2").Value = 94
$ws.Range("I7").Value = 281

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("I2").Value = 84
$ws.Range("H3").Value = 336
$ws.Range("I3").Value = 79
$ws.Range("H7").Value = 941
$ws.Range("I7").Value = 267

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("I3").Value = 117
$ws.Range("I6").Value = 110
$ws.Range("I7").Value = 325

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("I2").Value = 25
$ws.Range("I3").Value = 24
$ws.Range("I7").Value = 82

$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("I2").Value = 24
$ws.Range("I7").Value = 62

$ws = $wb.Worksheets.Item("New City")
$ws.Range("I2").Value = 60
$ws.Range("I5").Value = 10
$ws.Range("I7").Value = 198

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("I2").Value = 57
$ws.Range("I3").Value = 66
$ws.Range("I6").Value = 27
$ws.Range("I7").Value = 164

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("I2").Value = 44
$ws.Range("I3").Value = 56
$ws.Range("I7").Value = 140

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("I2").Value = 92
$ws.Range("I3").Value = 133
$ws.Range("I6").Value = 131
$ws.Range("I7").Value = 383

$ws = $wb.Worksheets.Item("Lincoln Park")
$ws.Range("I2").Value = 13
$ws.Range("I6").Value = 32
$ws.Range("I7").Value = 57

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("I2").Value = 43
$ws.Range("I7").Value = 189

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("I2").Value = 170
$ws.Range("I3").Value = 185
$ws.Range("I4").Value = 17
$ws.Range("I7").Value = 542

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("I3").Value = 62
$ws.Range("I6").Value = 67
$ws.Range("I7").Value = 236

$ws = $wb.Worksheets.Item("Irving Park")
$ws.Range("I3").Value = 15
$ws.Range("I7").Value = 66

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("I6").Value = 53
$ws.Range("I7").Value = 91

$ws = $wb.Worksheets.Item("River North")
$ws.Range("I2").Value = 23
$ws.Range("I7").Value = 129

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("I3").Value = 102
$ws.Range("I6").Value = 76
$ws.Range("I7").Value = 291

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("I6").Value = 29
$ws.Range("I7").Value = 97

$ws = $wb.Worksheets.Item("Chinatown")
$ws.Range("I6").Value = 42
$ws.Range("I7").Value = 51

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("I3").Value = 66
$ws.Range("I6").Value = 75
$ws.Range("I7").Value = 219

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("G4").Value = 25
$ws.Range("I5").Value = 7
$ws.Range("G7").Value = 619
$ws.Range("I7").Value = 222

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("I2").Value = 34
$ws.Range("I6").Value = 34
$ws.Range("I7").Value = 113

$ws = $wb.Worksheets.Item("Wrigleyville")
$ws.Range("I2").Value = 2
$ws.Range("I6").Value = 13

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("I2").Value = 15
$ws.Range("I7").Value = 73

$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Range("I3").Value = 27
$ws.Range("I7").Value = 109

$ws = $wb.Worksheets.Item("Avalon Park")
$ws.Range("I2").Value = 17
$ws.Range("I7").Value = 45

$ws = $wb.Worksheets.Item("Streeterville")
$ws.Range("I5").Value = 2
$ws.Range("I7").Value = 49

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("I2").Value = 31
$ws.Range("I6").Value = 37
$ws.Range("I7").Value = 94

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("I2").Value = 13
$ws.Range("I6").Value = 37
$ws.Range("I7").Value = 76

$ws = $wb.Worksheets.Item("Riverdale")
$ws.Range("I6").Value = 11
$ws.Range("I7").Value = 43

$ws = $wb.Worksheets.Item("Ukrainian Village")
$ws.Range("I6").Value = 6
$ws.Range("I7").Value = 15
